# Refresh the cryptos list with updated Price / Volume(1h) figures (and
# restore the correct Coin/Link pairing for two rows that had swapped),
# mirroring the periodic "Updated cryptos list ... with GitHub Actions" run.
#
# Numeric-looking Price values are written with a leading apostrophe so
# Excel keeps them as literal text (matching the source sheet, which
# stores every Price/Volume cell as a string, e.g. "3.800.67" style
# thousand-separated prices that are not valid numbers anyway).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.913.99'
$ws.Range("E2").Value = '  -2.13%  '
$ws.Range("D3").Value = '3.803.08'
$ws.Range("E3").Value = '  +0.90%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'601.52"
$ws.Range("E5").Value = '  -2.08%  '
$ws.Range("D6").Value = "'172.30"
$ws.Range("E6").Value = '  -3.34%  '
$ws.Range("D7").Value = '3.800.78'
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("E9").Value = '  +0.47%  '
$ws.Range("E10").Value = '  -4.04%  '
$ws.Range("D11").Value = "'6.20"
$ws.Range("E11").Value = '  -6.11%  '
$ws.Range("E12").Value = '  -3.69%  '
$ws.Range("D13").Value = "'38.83"
$ws.Range("E13").Value = '  -2.72%  '
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("D15").Value = '4.436.82'
$ws.Range("E15").Value = '  +0.95%  '
$ws.Range("D16").Value = '3.801.68'
$ws.Range("E16").Value = '  +0.90%  '
$ws.Range("D17").Value = '67.882.04'
$ws.Range("E17").Value = '  -2.27%  '
$ws.Range("E18").Value = '  -3.61%  '
$ws.Range("E19").Value = '  -3.93%  '
$ws.Range("D20").Value = "'17.28"
$ws.Range("E20").Value = '  +5.72%  '
$ws.Range("D21").Value = "'494.16"
$ws.Range("E21").Value = '  -2.67%  '
$ws.Range("D22").Value = "'9.21"
$ws.Range("E22").Value = '  -2.18%  '
$ws.Range("D23").Value = "'0.743"
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").Value = "'85.78"
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("E25").Value = '  -4.49%  '
$ws.Range("E26").Value = '  +7.37%  '
$ws.Range("E27").Value = '  -3.57%  '
$ws.Range("D28").Value = "'10.22"
$ws.Range("E28").Value = '  -3.69%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  -0.37%  '
$ws.Range("D31").Value = "'2.45"
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("D32").Value = "'32.98"
$ws.Range("E32").Value = '  +7.53%  '
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("E34").Value = '  -3.64%  '
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").Value = "'1.02"
$ws.Range("E36").Value = '  -3.14%  '
$ws.Range("E37").Value = '  -4.56%  '
$ws.Range("B38").Value = 'TheGraph'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D38").Value = "'0.332"
$ws.Range("E38").Value = '  -2.43%  '
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = "'463.57"
$ws.Range("E39").Value = '  +1.35%  '
$ws.Range("E40").Value = '  -5.60%  '
$ws.Range("E41").Value = '  -2.24%  '
$ws.Range("D42").Value = "'49.05"
$ws.Range("E42").Value = '  -1.49%  '
$ws.Range("D43").Value = "'2.86"
$ws.Range("E43").Value = '  -3.61%  '
$ws.Range("D44").Value = "'8.45"
$ws.Range("D45").Value = "'41.11"
$ws.Range("E45").Value = '  -7.72%  '
$ws.Range("D47").Value = '2.845.54'
$ws.Range("E47").Value = '  -3.86%  '
$ws.Range("D48").Value = "'139.73"
$ws.Range("E48").Value = '  +0.54%  '
$ws.Range("D49").Value = "'0.0352"
$ws.Range("E49").Value = '  -2.28%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'24.48"
$ws.Range("E50").Value = '  +12.93%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = "'25.91"
$ws.Range("E51").Value = '  -5.16%  '
